$wb = $excel.ActiveWorkbook

$website    = $wb.Worksheets.Item("Website")
$everything = $wb.Worksheets.Item("Everything")

# --- "Website" sheet, row 3: I Am Entertainment Magazine ---
$website.Range("A3").Value = "About I Am Entertainment Magazine"
$website.Range("B3").Value = "info@iaemagazine.com"
$website.Range("C3").Value = "http://www.iaemagazine.com"
$website.Range("D3").Value = "United States"
$website.Range("E3").Value = "United States"
$website.Range("F3").Value = "N/A"
$website.Range("G3").Value = "N/A"
$website.Range("H3").Value = "N/A"
$website.Range("I3").Value = "N/A"
$website.Range("J3").Value = "N/A"
$website.Range("K3").Value = "N/A"
$website.Range("L3").Value = "Check the website"

# --- "Everything" sheet, row 5: Alarm Magazine ---
$everything.Range("A5").Value = "Alarm Magazine"
$everything.Range("B5").Value = "scottm@alarmpress.com"
$everything.Range("C5").Value = "http://alarm-magazine.com"
$everything.Range("D5").Value = "North America"
$everything.Range("E5").Value = "United States"
$everything.Range("F5").Value = "N/A"
$everything.Range("G5").Value = "N/A"
$everything.Range("H5").Value = "N/A"
$everything.Range("I5").Value = "Attn: Music Editor, 900 N. Franklin St. #300, Chicago, IL 60610"
$everything.Range("J5").Value = 3123411301
$everything.Range("K5").Value = "N/A"
$everything.Range("L5").Value = "They want to know about shows. "

# --- "Website" sheet, row 5: The A.V. Club ---
$website.Range("A5").Value = "The A.V. Club"
$website.Range("B5").Value = "N/A"
$website.Range("C5").Value = "http://www.avclub.com"
$website.Range("D5").Value = "North America"
$website.Range("E5").Value = "United States"
$website.Range("F5").Value = "N/A"
$website.Range("G5").Value = "N/A"
$website.Range("H5").Value = "N/A"
$website.Range("I5").Value = "N/A"
$website.Range("J5").Value = "N/A"
$website.Range("K5").Value = "N/A"
$website.Range("L5").Value = "Interviews, essays and reviews of movies, music and books."

# The workbook was last saved with the "Website" tab active/selected
# (previously it was "Mp3").
$website.Activate()
